$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$newText = "The project aims to address and analyse behavioral data, where we can predict treatment outcomes and potential risks for chronic illness. The health data collected can be used for risk scoring, readmission prediction and prevention, predicting infection and deterioration and so much more at the individual patient level."

# Replace the whole paragraph text (keeps the existing run formatting: italic,
# size 2000, tx1 scheme color w/ lumMod/lumOff).
$tr.Text = $newText

# Split the single run into three runs so that "behavioral" becomes its own
# run (mirrors the authored edit marking that word as a flagged spelling
# variant). Re-asserting the already-true Italic formatting on the substring
# forces PowerPoint to break the run without changing any visible formatting.
$start = $newText.IndexOf("behavioral") + 1
$len = "behavioral".Length
$mid = $tr.Characters($start, $len)
$mid.Font.Italic = $true
